{"js": "// Replace the 25 \"three-digit x one-digit\" practice equations in the\n// single table with the newly generated set, cell by cell (row, col),\n// matching document order exactly as in the diff.\nconst table = context.document.body.tables.getFirst();\n\n// Rows that actually carry equations (every 5th row, the others are blank\n// spacer rows) mapped to the 5 new values for that row, left to right.\nconst updates = [\n  { row: 0, values: [\"148\u00d77=1036\", \"210\u00d77=1470\", \"189\u00d75=945\", \"887\u00d73=2661\", \"795\u00d73=2385\"] },\n  { row: 4, values: [\"983\u00d77=6881\", \"710\u00d78=5680\", \"292\u00d74=1168\", \"231\u00d72=462\", \"735\u00d76=4410\"] },\n  { row: 9, values: [\"472\u00d79=4248\", \"495\u00d74=1980\", \"666\u00d78=5328\", \"374\u00d75=1870\", \"250\u00d76=1500\"] },\n  { row: 14, values: [\"874\u00d73=2622\", \"793\u00d74=3172\", \"821\u00d72=1642\", \"185\u00d78=1480\", \"135\u00d75=675\"] },\n  { row: 19, values: [\"904\u00d72=1808\", \"283\u00d73=849\", \"310\u00d76=1860\", \"246\u00d77=1722\", \"156\u00d75=780\"] },\n];\n\nfor (const { row, values } of updates) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" practice equations in the\n# single table with the newly generated set, cell by cell (1-based\n# Table.Cell(row, col) addressing), matching document order exactly as\n# in the diff. Blank spacer rows (2-4, 6-9, 11-14, 16-19) are untouched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Values = @(\"148\u00d77=1036\", \"210\u00d77=1470\", \"189\u00d75=945\", \"887\u00d73=2661\", \"795\u00d73=2385\") },\n    @{ Row = 5;  Values = @(\"983\u00d77=6881\", \"710\u00d78=5680\", \"292\u00d74=1168\", \"231\u00d72=462\",  \"735\u00d76=4410\") },\n    @{ Row = 10; Values = @(\"472\u00d79=4248\", \"495\u00d74=1980\", \"666\u00d78=5328\", \"374\u00d75=1870\", \"250\u00d76=1500\") },\n    @{ Row = 15; Values = @(\"874\u00d73=2622\", \"793\u00d74=3172\", \"821\u00d72=1642\", \"185\u00d78=1480\", \"135\u00d75=675\")  },\n    @{ Row = 20; Values = @(\"904\u00d72=1808\", \"283\u00d73=849\",  \"310\u00d76=1860\", \"246\u00d77=1722\", \"156\u00d75=780\")  }\n)\n\nforeach ($u in $updates) {\n    $row = $u.Row\n    $vals = $u.Values\n    for ($col = 1; $col -le $vals.Count; $col++) {\n        $t.Cell($row, $col).Range.Text = $vals[$col - 1]\n    }\n}\n"}
